$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 796.1818
$ws.Cells.Item(6, 9).Value = 189.75
$ws.Cells.Item(6, 10).Value = 2413.3333
$ws.Cells.Item(6, 11).Value = 569.25
$ws.Cells.Item(6, 12).Value = 7239.999899999999
$ws.Cells.Item(6, 13).Value = -457.25
$ws.Cells.Item(6, 14).Value = -7463.999899999999

$ws.Cells.Item(12, 8).Value = 1163.8334
$ws.Cells.Item(12, 9).Value = 1030.6666
$ws.Cells.Item(12, 10).Value = 1563.3334
$ws.Cells.Item(12, 11).Value = 1030.6666
$ws.Cells.Item(12, 12).Value = 1563.3334
$ws.Cells.Item(12, 13).Value = -860.6666
$ws.Cells.Item(12, 14).Value = -1903.3334

$ws.Cells.Item(51, 8).Value = 6224.4707
$ws.Cells.Item(51, 10).Value = 6910
$ws.Cells.Item(51, 12).Value = 6910
$ws.Cells.Item(51, 14).Value = -7878

$ws.Cells.Item(57, 8).Value = 139887.17
$ws.Cells.Item(57, 10).Value = 139887.17
$ws.Cells.Item(57, 12).Value = 419661.51
$ws.Cells.Item(57, 14).Value = -420659.51

$ws.Cells.Item(62, 8).Value = 4412.625
$ws.Cells.Item(62, 9).Value = 3349.4
$ws.Cells.Item(62, 11).Value = 3349.4
$ws.Cells.Item(62, 13).Value = -2725.4

$ws.Cells.Item(65, 8).Value = 4412.625
$ws.Cells.Item(65, 9).Value = 3349.4
$ws.Cells.Item(65, 11).Value = 16747
$ws.Cells.Item(65, 13).Value = -13627

$ws.Cells.Item(76, 8).Value = 14459.083
$ws.Cells.Item(76, 10).Value = 12699.857
$ws.Cells.Item(76, 12).Value = 12699.857
$ws.Cells.Item(76, 14).Value = -13329.857

$ws.Cells.Item(79, 8).Value = 14459.083
$ws.Cells.Item(79, 10).Value = 12699.857
$ws.Cells.Item(79, 12).Value = 12699.857
$ws.Cells.Item(79, 14).Value = -14883.857

$ws.Cells.Item(112, 8).Value = 8000
$ws.Cells.Item(112, 9).Value = 8000
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 24000
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = -22892
$ws.Cells.Item(112, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 14337.091
$ws.Cells.Item(113, 10).Value = 28124.75
$ws.Cells.Item(113, 12).Value = 28124.75
$ws.Cells.Item(113, 14).Value = -34632.75

$ws.Cells.Item(128, 8).Value = 174987
$ws.Cells.Item(128, 10).Value = 174987
$ws.Cells.Item(128, 12).Value = 174987
$ws.Cells.Item(128, 14).Value = -184947

$ws.Cells.Item(136, 8).Value = 179994.8
$ws.Cells.Item(136, 10).Value = 179994.8
$ws.Cells.Item(136, 12).Value = 179994.8
$ws.Cells.Item(136, 14).Value = -190194.8

$ws.Cells.Item(137, 8).Value = 7435.4443
$ws.Cells.Item(137, 9).Value = 7577.905
$ws.Cells.Item(137, 10).Value = 7236
$ws.Cells.Item(137, 11).Value = 22733.715
$ws.Cells.Item(137, 12).Value = 21708
$ws.Cells.Item(137, 13).Value = -20183.715
$ws.Cells.Item(137, 14).Value = -26808

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 52859004
$ws.Cells.Item(2, 9).Value = 58731950
$ws.Cells.Item(2, 10).Value = 2500
$ws.Cells.Item(2, 11).Value = 58731950
$ws.Cells.Item(2, 12).Value = 2500
$ws.Cells.Item(2, 13).Value = -58731837
$ws.Cells.Item(2, 14).Value = -2726

$ws.Cells.Item(32, 8).Value = 4100.195
$ws.Cells.Item(32, 9).Value = 2996.8242
$ws.Cells.Item(32, 10).Value = 31316.666
$ws.Cells.Item(32, 11).Value = 2996.8242
$ws.Cells.Item(32, 12).Value = 31316.666
$ws.Cells.Item(32, 13).Value = -2709.8242
$ws.Cells.Item(32, 14).Value = -31890.666

$ws.Cells.Item(45, 8).Value = 600
$ws.Cells.Item(45, 9).Value = 600
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 600
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -223
$ws.Cells.Item(45, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 15320.929
$ws.Cells.Item(61, 9).Value = 13449.8
$ws.Cells.Item(61, 10).Value = 19998.75
$ws.Cells.Item(61, 11).Value = 13449.8
$ws.Cells.Item(61, 12).Value = 19998.75
$ws.Cells.Item(61, 13).Value = -13237.8
$ws.Cells.Item(61, 14).Value = -20422.75

$ws.Cells.Item(63, 8).Value = 3350.8125
$ws.Cells.Item(63, 10).Value = 3517.1667
$ws.Cells.Item(63, 12).Value = 3517.1667
$ws.Cells.Item(63, 14).Value = -4889.1667

$ws.Cells.Item(66, 8).Value = 3350.8125
$ws.Cells.Item(66, 10).Value = 3517.1667
$ws.Cells.Item(66, 12).Value = 17585.8335
$ws.Cells.Item(66, 14).Value = -24449.8335

$ws.Cells.Item(88, 8).Value = 1489.6
$ws.Cells.Item(88, 10).Value = 1655.3334
$ws.Cells.Item(88, 12).Value = 1655.3334
$ws.Cells.Item(88, 14).Value = -2467.3334

$ws.Cells.Item(91, 8).Value = 1489.6
$ws.Cells.Item(91, 10).Value = 1655.3334
$ws.Cells.Item(91, 12).Value = 1655.3334
$ws.Cells.Item(91, 14).Value = -4463.3334

$ws.Cells.Item(102, 8).Value = 2827.25
$ws.Cells.Item(102, 9).Value = 2818.158
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 2818.158
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -1196.158
$ws.Cells.Item(102, 14).Value = -6244

$ws.Cells.Item(116, 8).Value = 52859004
$ws.Cells.Item(116, 9).Value = 58731950
$ws.Cells.Item(116, 10).Value = 2500
$ws.Cells.Item(116, 11).Value = 58731950
$ws.Cells.Item(116, 12).Value = 2500
$ws.Cells.Item(116, 13).Value = -58729656
$ws.Cells.Item(116, 14).Value = -7088

$ws.Cells.Item(118, 8).Value = 40000
$ws.Cells.Item(118, 10).Value = 40000
$ws.Cells.Item(118, 12).Value = 40000
$ws.Cells.Item(118, 14).Value = -43314

$ws.Cells.Item(122, 8).Value = 1996.8064
$ws.Cells.Item(122, 9).Value = 1472.5652
$ws.Cells.Item(122, 11).Value = 4417.6956
$ws.Cells.Item(122, 13).Value = -1967.6956

$ws.Cells.Item(132, 8).Value = 4497.5835
$ws.Cells.Item(132, 9).Value = 4197.2
$ws.Cells.Item(132, 10).Value = 5999.5
$ws.Cells.Item(132, 11).Value = 12591.6
$ws.Cells.Item(132, 12).Value = 17998.5
$ws.Cells.Item(132, 13).Value = -10061.6
$ws.Cells.Item(132, 14).Value = -23058.5

$ws.Cells.Item(136, 8).Value = 15320.929
$ws.Cells.Item(136, 9).Value = 13449.8
$ws.Cells.Item(136, 10).Value = 19998.75
$ws.Cells.Item(136, 11).Value = 40349.39999999999
$ws.Cells.Item(136, 12).Value = 59996.25
$ws.Cells.Item(136, 13).Value = -37799.39999999999
$ws.Cells.Item(136, 14).Value = -65096.25

$ws.Cells.Item(141, 8).Value = 79160.375
$ws.Cells.Item(141, 10).Value = 79160.375
$ws.Cells.Item(141, 12).Value = 79160.375
$ws.Cells.Item(141, 14).Value = -89520.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 52859004
$ws.Cells.Item(3, 9).Value = 58731950
$ws.Cells.Item(3, 10).Value = 2500
$ws.Cells.Item(3, 11).Value = 58731950
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = -58731836
$ws.Cells.Item(3, 14).Value = -2728

$ws.Cells.Item(86, 8).Value = 2083.742
$ws.Cells.Item(86, 9).Value = 2161.682
$ws.Cells.Item(86, 10).Value = 1893.2222
$ws.Cells.Item(86, 11).Value = 2161.682
$ws.Cells.Item(86, 12).Value = 1893.2222
$ws.Cells.Item(86, 13).Value = -1038.682
$ws.Cells.Item(86, 14).Value = -4139.2222

$ws.Cells.Item(89, 8).Value = 2083.742
$ws.Cells.Item(89, 9).Value = 2161.682
$ws.Cells.Item(89, 10).Value = 1893.2222
$ws.Cells.Item(89, 11).Value = 10808.41
$ws.Cells.Item(89, 12).Value = 9466.110999999999
$ws.Cells.Item(89, 13).Value = -5192.41
$ws.Cells.Item(89, 14).Value = -20698.111

$ws.Cells.Item(119, 8).Value = 93799.60000000001
$ws.Cells.Item(119, 10).Value = 93799.60000000001
$ws.Cells.Item(119, 12).Value = 93799.60000000001
$ws.Cells.Item(119, 14).Value = -103475.6

$ws.Cells.Item(132, 8).Value = 113360
$ws.Cells.Item(132, 10).Value = 118840
$ws.Cells.Item(132, 12).Value = 118840
$ws.Cells.Item(132, 14).Value = -128960

$ws.Cells.Item(134, 8).Value = 2551.2322
$ws.Cells.Item(134, 9).Value = 2586.7273
$ws.Cells.Item(134, 11).Value = 7760.1819
$ws.Cells.Item(134, 13).Value = -5225.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 293.78946
$ws.Cells.Item(7, 9).Value = 16.538462
$ws.Cells.Item(7, 10).Value = 894.5
$ws.Cells.Item(7, 11).Value = 16.538462
$ws.Cells.Item(7, 12).Value = 894.5
$ws.Cells.Item(7, 13).Value = 96.461538
$ws.Cells.Item(7, 14).Value = -1120.5

$ws.Cells.Item(31, 8).Value = 2056.6667
$ws.Cells.Item(31, 9).Value = 1237.5
$ws.Cells.Item(31, 11).Value = 1237.5
$ws.Cells.Item(31, 13).Value = -942.5

$ws.Cells.Item(34, 8).Value = 2056.6667
$ws.Cells.Item(34, 9).Value = 1237.5
$ws.Cells.Item(34, 11).Value = 1237.5
$ws.Cells.Item(34, 13).Value = -1035.5

$ws.Cells.Item(58, 8).Value = 8822.695
$ws.Cells.Item(58, 9).Value = 7909.4287
$ws.Cells.Item(58, 10).Value = 10243.333
$ws.Cells.Item(58, 11).Value = 7909.4287
$ws.Cells.Item(58, 12).Value = 10243.333
$ws.Cells.Item(58, 13).Value = -7706.4287
$ws.Cells.Item(58, 14).Value = -10649.333

$ws.Cells.Item(110, 8).Value = 99999
$ws.Cells.Item(110, 10).Value = 99999
$ws.Cells.Item(110, 12).Value = 99999
$ws.Cells.Item(110, 14).Value = -108179

$ws.Cells.Item(135, 8).Value = 89399.39999999999
$ws.Cells.Item(135, 10).Value = 89399.39999999999
$ws.Cells.Item(135, 12).Value = 89399.39999999999
$ws.Cells.Item(135, 14).Value = -99539.39999999999

$ws.Cells.Item(136, 8).Value = 8822.695
$ws.Cells.Item(136, 9).Value = 7909.4287
$ws.Cells.Item(136, 10).Value = 10243.333
$ws.Cells.Item(136, 11).Value = 23728.2861
$ws.Cells.Item(136, 12).Value = 30729.999
$ws.Cells.Item(136, 13).Value = -21178.2861
$ws.Cells.Item(136, 14).Value = -35829.999

$ws.Cells.Item(140, 8).Value = 114526.375
$ws.Cells.Item(140, 10).Value = 114495.14
$ws.Cells.Item(140, 12).Value = 114495.14
$ws.Cells.Item(140, 14).Value = -124855.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 43506.75
$ws.Cells.Item(4, 10).Value = 3264.3333
$ws.Cells.Item(4, 12).Value = 9792.999899999999
$ws.Cells.Item(4, 14).Value = -10016.9999

$ws.Cells.Item(7, 8).Value = 130
$ws.Cells.Item(7, 9).Value = 56
$ws.Cells.Item(7, 11).Value = 168
$ws.Cells.Item(7, 13).Value = -56

$ws.Cells.Item(11, 8).Value = 1576.6538
$ws.Cells.Item(11, 9).Value = 587.3333
$ws.Cells.Item(11, 11).Value = 1761.9999
$ws.Cells.Item(11, 13).Value = -1621.9999

$ws.Cells.Item(12, 8).Value = 221.625
$ws.Cells.Item(12, 10).Value = 241.5
$ws.Cells.Item(12, 12).Value = 724.5
$ws.Cells.Item(12, 14).Value = -1070.5

$ws.Cells.Item(62, 8).Value = 9933.556
$ws.Cells.Item(62, 9).Value = 4756
$ws.Cells.Item(62, 10).Value = 11412.857
$ws.Cells.Item(62, 11).Value = 14268
$ws.Cells.Item(62, 12).Value = 34238.571
$ws.Cells.Item(62, 13).Value = -13582
$ws.Cells.Item(62, 14).Value = -35610.571

$ws.Cells.Item(65, 8).Value = 9933.556
$ws.Cells.Item(65, 9).Value = 4756
$ws.Cells.Item(65, 10).Value = 11412.857
$ws.Cells.Item(65, 11).Value = 42804
$ws.Cells.Item(65, 12).Value = 102715.713
$ws.Cells.Item(65, 13).Value = -39372
$ws.Cells.Item(65, 14).Value = -109579.713

$ws.Cells.Item(103, 8).Value = 274.08334
$ws.Cells.Item(103, 9).Value = 74.40000000000001
$ws.Cells.Item(103, 10).Value = 416.7143
$ws.Cells.Item(103, 11).Value = 223.2
$ws.Cells.Item(103, 12).Value = 1250.1429
$ws.Cells.Item(103, 13).Value = 655.8
$ws.Cells.Item(103, 14).Value = -3008.1429

$ws.Cells.Item(131, 8).Value = 1015593.1
$ws.Cells.Item(131, 9).Value = 1839320
$ws.Cells.Item(131, 10).Value = 1775.3077
$ws.Cells.Item(131, 11).Value = 5517960
$ws.Cells.Item(131, 12).Value = 5325.9231
$ws.Cells.Item(131, 13).Value = -5512920
$ws.Cells.Item(131, 14).Value = -15405.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(58, 8).Value = 58375.8
$ws.Cells.Item(58, 10).Value = 58375.8
$ws.Cells.Item(58, 12).Value = 58375.8
$ws.Cells.Item(58, 14).Value = -58929.8

$ws.Cells.Item(64, 8).Value = 60000
$ws.Cells.Item(64, 10).Value = 60000
$ws.Cells.Item(64, 12).Value = 60000
$ws.Cells.Item(64, 14).Value = -60496

$ws.Cells.Item(67, 8).Value = 60000
$ws.Cells.Item(67, 10).Value = 60000
$ws.Cells.Item(67, 12).Value = 60000
$ws.Cells.Item(67, 14).Value = -61716

$ws.Cells.Item(70, 8).Value = 5042.364
$ws.Cells.Item(70, 9).Value = 4261.6665
$ws.Cells.Item(70, 10).Value = 5979.2
$ws.Cells.Item(70, 11).Value = 4261.6665
$ws.Cells.Item(70, 12).Value = 5979.2
$ws.Cells.Item(70, 13).Value = -3991.6665
$ws.Cells.Item(70, 14).Value = -6519.2

$ws.Cells.Item(73, 8).Value = 5042.364
$ws.Cells.Item(73, 9).Value = 4261.6665
$ws.Cells.Item(73, 10).Value = 5979.2
$ws.Cells.Item(73, 11).Value = 4261.6665
$ws.Cells.Item(73, 12).Value = 5979.2
$ws.Cells.Item(73, 13).Value = -3325.6665
$ws.Cells.Item(73, 14).Value = -7851.2

$ws.Cells.Item(80, 8).Value = 44377332
$ws.Cells.Item(80, 10).Value = 2569.6667
$ws.Cells.Item(80, 12).Value = 2569.6667
$ws.Cells.Item(80, 14).Value = -4565.6667

$ws.Cells.Item(83, 8).Value = 44377332
$ws.Cells.Item(83, 10).Value = 2569.6667
$ws.Cells.Item(83, 12).Value = 12848.3335
$ws.Cells.Item(83, 14).Value = -22832.3335

$ws.Cells.Item(99, 8).Value = 11058.5
$ws.Cells.Item(99, 10).Value = 17999
$ws.Cells.Item(99, 12).Value = 17999
$ws.Cells.Item(99, 14).Value = -22491

$ws.Cells.Item(102, 8).Value = 7298.3076
$ws.Cells.Item(102, 9).Value = 7426.8
$ws.Cells.Item(102, 11).Value = 7426.8
$ws.Cells.Item(102, 13).Value = -5804.8

$ws.Cells.Item(122, 8).Value = 3789.3794
$ws.Cells.Item(122, 9).Value = 3674.7144
$ws.Cells.Item(122, 10).Value = 7000
$ws.Cells.Item(122, 11).Value = 11024.1432
$ws.Cells.Item(122, 12).Value = 21000
$ws.Cells.Item(122, 13).Value = -8574.143199999999
$ws.Cells.Item(122, 14).Value = -25900

$ws.Cells.Item(126, 8).Value = 4191.951
$ws.Cells.Item(126, 10).Value = 3677.6365
$ws.Cells.Item(126, 12).Value = 11032.9095
$ws.Cells.Item(126, 14).Value = -15972.9095

$ws.Cells.Item(132, 8).Value = 2465.6667
$ws.Cells.Item(132, 9).Value = 2465.6667
$ws.Cells.Item(132, 11).Value = 7397.000100000001
$ws.Cells.Item(132, 13).Value = -4867.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 13893850
$ws.Cells.Item(16, 9).Value = 15629344
$ws.Cells.Item(16, 10).Value = 9899
$ws.Cells.Item(16, 11).Value = 15629344
$ws.Cells.Item(16, 12).Value = 9899
$ws.Cells.Item(16, 13).Value = -15629174
$ws.Cells.Item(16, 14).Value = -10239

$ws.Cells.Item(55, 8).Value = 2493.75
$ws.Cells.Item(55, 9).Value = 2357.5
$ws.Cells.Item(55, 10).Value = 2698.125
$ws.Cells.Item(55, 11).Value = 2357.5
$ws.Cells.Item(55, 12).Value = 2698.125
$ws.Cells.Item(55, 13).Value = -2184.5
$ws.Cells.Item(55, 14).Value = -3044.125

$ws.Cells.Item(61, 8).Value = 10886.333
$ws.Cells.Item(61, 9).Value = 11521.454
$ws.Cells.Item(61, 10).Value = 3900
$ws.Cells.Item(61, 11).Value = 11521.454
$ws.Cells.Item(61, 12).Value = 3900
$ws.Cells.Item(61, 13).Value = -11319.454
$ws.Cells.Item(61, 14).Value = -4304

$ws.Cells.Item(113, 8).Value = 10886.333
$ws.Cells.Item(113, 9).Value = 11521.454
$ws.Cells.Item(113, 10).Value = 3900
$ws.Cells.Item(113, 11).Value = 11521.454
$ws.Cells.Item(113, 12).Value = 3900
$ws.Cells.Item(113, 13).Value = -9351.454
$ws.Cells.Item(113, 14).Value = -8240

$ws.Cells.Item(132, 8).Value = 33183.973
$ws.Cells.Item(132, 9).Value = 52021.5
$ws.Cells.Item(132, 10).Value = 3582.1428
$ws.Cells.Item(132, 11).Value = 156064.5
$ws.Cells.Item(132, 12).Value = 10746.4284
$ws.Cells.Item(132, 13).Value = -153534.5
$ws.Cells.Item(132, 14).Value = -15806.4284

$ws.Cells.Item(133, 8).Value = 120000
$ws.Cells.Item(133, 10).Value = 120000
$ws.Cells.Item(133, 12).Value = 120000
$ws.Cells.Item(133, 14).Value = -125060

$ws.Cells.Item(136, 8).Value = 5640646
$ws.Cells.Item(136, 9).Value = 6219030.5
$ws.Cells.Item(136, 10).Value = 49595.668
$ws.Cells.Item(136, 11).Value = 18657091.5
$ws.Cells.Item(136, 12).Value = 148787.004
$ws.Cells.Item(136, 13).Value = -18654541.5
$ws.Cells.Item(136, 14).Value = -153887.004

$ws.Cells.Item(137, 8).Value = 126601.2
$ws.Cells.Item(137, 10).Value = 108251.75
$ws.Cells.Item(137, 12).Value = 108251.75
$ws.Cells.Item(137, 14).Value = -118451.75

$ws.Cells.Item(138, 8).Value = 129455.29
$ws.Cells.Item(138, 10).Value = 129455.29
$ws.Cells.Item(138, 12).Value = 129455.29
$ws.Cells.Item(138, 14).Value = -139735.29

$ws.Cells.Item(139, 8).Value = 45124.75
$ws.Cells.Item(139, 10).Value = 45124.75
$ws.Cells.Item(139, 12).Value = 45124.75
$ws.Cells.Item(139, 14).Value = -55404.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2007
$ws.Cells.Item(132, 9).Value = 1586.2354
$ws.Cells.Item(132, 10).Value = 7371.75
$ws.Cells.Item(132, 11).Value = 4758.706200000001
$ws.Cells.Item(132, 12).Value = 22115.25
$ws.Cells.Item(132, 13).Value = -2228.706200000001
$ws.Cells.Item(132, 14).Value = -27175.25

$ws.Cells.Item(136, 8).Value = 1940.4
$ws.Cells.Item(136, 9).Value = 1983.4572
$ws.Cells.Item(136, 11).Value = 5950.3716
$ws.Cells.Item(136, 13).Value = -3400.3716

$ws.Cells.Item(137, 8).Value = 98927.25
$ws.Cells.Item(137, 10).Value = 98927.25
$ws.Cells.Item(137, 12).Value = 98927.25
$ws.Cells.Item(137, 14).Value = -109127.25
